$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = "BIM"
$ws.Range("B4").Value = 147
$ws.Range("C4").Value = 96
$ws.Range("D4").Value = 4
$ws.Range("E4").Value = 20200504
$ws.Range("F4").Value = "3cm branch has died"

$ws.Range("B5").Select()
